{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: In the \"Crafting\" intro paragraph, fix the wording from\n//   \"...Cost and the Quality, but the Time will be the variable...\"\n// to\n//   \"...Cost and the Quality, and Time will be the variable...\"\n// (This also happens to merge the two runs the sentence was split across,\n// which is what Word naturally does when the replaced text spans a run\n// boundary.)\n//\n// Change 2: Make the section's page orientation explicitly Portrait\n// (it already renders as portrait - w=12240 < h=15840 - but the\n// orientation attribute is now written explicitly on pgSz).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst oldFragment = \"but the Time will be the\";\nconst newFragment = \"and Time will be the\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(oldFragment) !== -1) {\n    const newText = para.text.split(oldFragment).join(newFragment);\n    para.getRange().insertText(newText, Word.InsertLocation.replace);\n    break;\n  }\n}\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  sections.items[i].pageSetup.orientation = Word.PageOrientation.portrait;\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change 1: In the \"Crafting\" intro paragraph, fix the wording from\n#   \"...Cost and the Quality, but the Time will be the variable...\"\n# to\n#   \"...Cost and the Quality, and Time will be the variable...\"\n# (Replacing text that spans the original run boundary makes Word merge\n# the two runs into one, exactly like the authored edit.)\n#\n# Change 2: Make the section's page orientation explicitly Portrait\n# (already renders portrait, w=12240 < h=15840, but the orientation\n# attribute is now written explicitly on pgSz).\n\n$d = $word.ActiveDocument\n\n# Search text spans the original run boundary (the run break falls right\n# after \"...will be the \"), so Word collapses the paragraph back down to a\n# single run when it performs the replacement - matching the authored edit.\n$oldText = \"Usually, this will be the Cost and the Quality, but the Time will be the variable that successful Crafting Checks will determine, but you can shift this around if you want.\"\n$newText = \"Usually, this will be the Cost and the Quality, and Time will be the variable that successful Crafting Checks will determine, but you can shift this around if you want.\"\n\n$find = $d.Content.Find\n$find.Execute(\n    $oldText,\n    $false, $false, $false, $false, $false,\n    $true,\n    1,\n    $false,\n    $newText,\n    2\n)\n\nforeach ($sec in $d.Sections) {\n    $sec.PageSetup.Orientation = 0\n}\n"}
